# Apply attendance updates (săpt. 7 = column I, plus a handful of other
# catch-up checkmarks) to the "Prezenta_AnII_2024" sheet, add the new
# student "Alexandra Iovan" on row 57, and update the last active
# selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student row (57): name + săpt. 7 attendance
$ws.Range("B57").Value = "Alexandra Iovan"
$ws.Range("I57").Value = $true

# Row 4 (Alexandru Lupșe): săpt. 7 + săpt. 8
$ws.Range("I4").Value = $true
$ws.Range("J4").Value = $true

# Row 5 (Anamaria  Bohar): săpt. 7
$ws.Range("I5").Value = $true

# Row 7 (Anastasia Poleakova): săpt. 7
$ws.Range("I7").Value = $true

# Row 8 (Andreea Farcas): săpt. 5 + săpt. 7
$ws.Range("G8").Value = $true
$ws.Range("I8").Value = $true

# Row 9 (Andrei Mâța): săpt. 7
$ws.Range("I9").Value = $true

# Row 11 (Aniko Vieriu): săpt. 7
$ws.Range("I11").Value = $true

# Row 12 (Beniamin Vutan): săpt. 7
$ws.Range("I12").Value = $true

# Row 13 (Bianca Abrudan): săpt. 7
$ws.Range("I13").Value = $true

# Row 14 (Bianca Nicorici): săpt. 5 + săpt. 7
$ws.Range("G14").Value = $true
$ws.Range("I14").Value = $true

# Row 15 (Bogdan Bobos): săpt. 4 + săpt. 7
$ws.Range("F15").Value = $true
$ws.Range("I15").Value = $true

# Row 16 (Boglarka Szigeti): săpt. 7
$ws.Range("I16").Value = $true

# Row 17 (Casian Balaj): săpt. 2 + săpt. 7
$ws.Range("D17").Value = $true
$ws.Range("I17").Value = $true

# Row 19 (Cosmin Chira): săpt. 7
$ws.Range("I19").Value = $true

# Row 20 (Cristina Nemcea): săpt. 7
$ws.Range("I20").Value = $true

# Row 21 (Daniel Oistic): săpt. 7
$ws.Range("I21").Value = $true

# Row 22 (Daria Petre): săpt. 7
$ws.Range("I22").Value = $true

# Row 23 (Daria Puscas): săpt. 7
$ws.Range("I23").Value = $true

# Row 24 (Darius Orban): săpt. 5 + săpt. 7
$ws.Range("G24").Value = $true
$ws.Range("I24").Value = $true

# Row 26 (David Klein): săpt. 7
$ws.Range("I26").Value = $true

# Row 29 (Erik Lazin): săpt. 7
$ws.Range("I29").Value = $true

# Row 30 (Gabriel Talmazan): săpt. 7
$ws.Range("I30").Value = $true

# Row 32 (George Vezentan): săpt. 7
$ws.Range("I32").Value = $true

# Row 33 (Georgiana Galea): săpt. 7
$ws.Range("I33").Value = $true

# Row 34 (Kevin Csaba): săpt. 4 + săpt. 7
$ws.Range("F34").Value = $true
$ws.Range("I34").Value = $true

# Row 35 (Loris Cioban): săpt. 7
$ws.Range("I35").Value = $true

# Row 36 (Luca Bulea): săpt. 2 + săpt. 7
$ws.Range("D36").Value = $true
$ws.Range("I36").Value = $true

# Row 37 (Madalin Blaj): săpt. 7
$ws.Range("I37").Value = $true

# Row 40 (Noelia Sfrangeu): săpt. 2 + săpt. 7
$ws.Range("D40").Value = $true
$ws.Range("I40").Value = $true

# Row 42 (Patrick Tocut): săpt. 7
$ws.Range("I42").Value = $true

# Row 43 (Raluca Vereș): săpt. 7
$ws.Range("I43").Value = $true

# Row 45 (Raul Vonhala): săpt. 7
$ws.Range("I45").Value = $true

# Row 47 (Renata Halasz): săpt. 7
$ws.Range("I47").Value = $true

# Row 48 (Renata Tirban): săpt. 7
$ws.Range("I48").Value = $true

# Row 51 (Stefan Tulvan): săpt. 7
$ws.Range("I51").Value = $true

# Row 52 (Vanesa Clepce): săpt. 7
$ws.Range("I52").Value = $true

# Row 53 (Victor Pitirici): săpt. 7
$ws.Range("I53").Value = $true

# Row 54 (Vlad Brata): săpt. 7
$ws.Range("I54").Value = $true

# Update the last active selection shown in the bottom-right (frozen) pane
$ws.Range("K11").Select()
